$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (2-29) with corrected Beta calculations
# Row 2
$ws.Cells.Item(2, 3).Value = "Gau"
$ws.Cells.Item(2, 4).Value = 9.1868
$ws.Cells.Item(2, 5).Value = 68.5638
$ws.Cells.Item(2, 6).Value = 2.54
$ws.Cells.Item(2, 7).Value = 0.1339890729510325
$ws.Cells.Item(2, 8).Value = 48533979422205.02
$ws.Cells.Item(2, 9).Value = -1.6288

# Row 3
$ws.Cells.Item(3, 4).Value = 141.8609
$ws.Cells.Item(3, 5).Value = 163.0171
$ws.Cells.Item(3, 6).Value = 0.65
$ws.Cells.Item(3, 7).Value = 0.8702209768177693
$ws.Cells.Item(3, 8).Value = 10549768.4979
$ws.Cells.Item(3, 9).Value = 0.0125

# Row 4
$ws.Cells.Item(4, 4).Value = 180.5309
$ws.Cells.Item(4, 5).Value = 191.2037
$ws.Cells.Item(4, 6).Value = 0.97
$ws.Cells.Item(4, 7).Value = 0.9441809964974528
$ws.Cells.Item(4, 8).Value = 72932725.0883
$ws.Cells.Item(4, 9).Value = -0.0118

# Row 5
$ws.Cells.Item(5, 3).Value = "Exp"
$ws.Cells.Item(5, 4).Value = 97.7889
$ws.Cells.Item(5, 5).Value = 304.4531
$ws.Cells.Item(5, 6).Value = 0.45
$ws.Cells.Item(5, 7).Value = 0.3211952842654583
$ws.Cells.Item(5, 8).Value = 167601377.2439
$ws.Cells.Item(5, 9).Value = 0.0267

# Row 6
$ws.Cells.Item(6, 3).Value = "Sph"
$ws.Cells.Item(6, 4).Value = 188.0764
$ws.Cells.Item(6, 5).Value = 201.3225
$ws.Cells.Item(6, 6).Value = 1.82
$ws.Cells.Item(6, 7).Value = 0.9342045722658918
$ws.Cells.Item(6, 8).Value = 17580718.3172
$ws.Cells.Item(6, 9).Value = -0.0456

# Row 7
$ws.Cells.Item(7, 4).Value = 212.9603
$ws.Cells.Item(7, 5).Value = 212.9693
$ws.Cells.Item(7, 6).Value = 41.79
$ws.Cells.Item(7, 7).Value = 0.9999577403879336
$ws.Cells.Item(7, 8).Value = 81945553.5627
$ws.Cells.Item(7, 9).Value = -0.0469

# Row 8
$ws.Cells.Item(8, 3).Value = "Gau"
$ws.Cells.Item(8, 4).Value = 210.7451
$ws.Cells.Item(8, 5).Value = 2650116.7852
$ws.Cells.Item(8, 6).Value = 1966.72
$ws.Cells.Item(8, 7).Value = 0.00007952294826286134
$ws.Cells.Item(8, 8).Value = 52049794.887
$ws.Cells.Item(8, 9).Value = -7.4431

# Row 9
$ws.Cells.Item(9, 3).Value = "Sph"
$ws.Cells.Item(9, 4).Value = 149.0401
$ws.Cells.Item(9, 5).Value = 509.6292
$ws.Cells.Item(9, 6).Value = 55.3
$ws.Cells.Item(9, 7).Value = 0.2924481171800988
$ws.Cells.Item(9, 8).Value = 44208762.9972
$ws.Cells.Item(9, 9).Value = -1.7189

# Row 10
$ws.Cells.Item(10, 4).Value = 150.0045
$ws.Cells.Item(10, 5).Value = 208.4849
$ws.Cells.Item(10, 6).Value = 30.79
$ws.Cells.Item(10, 7).Value = 0.7194981507053988
$ws.Cells.Item(10, 8).Value = 2653164.836
$ws.Cells.Item(10, 9).Value = 0.391

# Row 12
$ws.Cells.Item(12, 4).Value = 67.7856
$ws.Cells.Item(12, 5).Value = 2051.8519
$ws.Cells.Item(12, 6).Value = 385.48
$ws.Cells.Item(12, 7).Value = 0.03303630247387738
$ws.Cells.Item(12, 8).Value = 17273879446.5675
$ws.Cells.Item(12, 9).Value = -0.8225

# Row 13
$ws.Cells.Item(13, 3).Value = "Exp"
$ws.Cells.Item(13, 4).Value = 167.2398
$ws.Cells.Item(13, 5).Value = 186.524
$ws.Cells.Item(13, 6).Value = 2.62
$ws.Cells.Item(13, 7).Value = 0.8966127683300809
$ws.Cells.Item(13, 8).Value = 323457.7288
$ws.Cells.Item(13, 9).Value = -0.015

# Row 14
$ws.Cells.Item(14, 4).Value = 259.5819
$ws.Cells.Item(14, 5).Value = 285.3169
$ws.Cells.Item(14, 6).Value = 2.09
$ws.Cells.Item(14, 7).Value = 0.9098020481787096
$ws.Cells.Item(14, 8).Value = 1290793.4812
$ws.Cells.Item(14, 9).Value = -0.153

# Row 15
$ws.Cells.Item(15, 4).Value = 285.1397
$ws.Cells.Item(15, 5).Value = 2781591.2417
$ws.Cells.Item(15, 6).Value = 933.61
$ws.Cells.Item(15, 7).Value = 0.0001025095620540327
$ws.Cells.Item(15, 8).Value = 111894461.2288
$ws.Cells.Item(15, 9).Value = -0.1606

# Row 16
$ws.Cells.Item(16, 3).Value = "Sph"
$ws.Cells.Item(16, 4).Value = 225.8824
$ws.Cells.Item(16, 5).Value = 250.564
$ws.Cells.Item(16, 6).Value = 8.210000000000001
$ws.Cells.Item(16, 7).Value = 0.9014958254178572
$ws.Cells.Item(16, 8).Value = 2775941.6876
$ws.Cells.Item(16, 9).Value = 0.3049

# Row 17
$ws.Cells.Item(17, 3).Value = "Sph"
$ws.Cells.Item(17, 4).Value = 95.06229999999999
$ws.Cells.Item(17, 5).Value = 918.0499
$ws.Cells.Item(17, 6).Value = 33.3
$ws.Cells.Item(17, 7).Value = 0.1035480751100784
$ws.Cells.Item(17, 8).Value = 24532592983.4695
$ws.Cells.Item(17, 9).Value = -7.0143

# Row 18
$ws.Cells.Item(18, 3).Value = "Sph"
$ws.Cells.Item(18, 4).Value = 284.4788
$ws.Cells.Item(18, 5).Value = 315.6029
$ws.Cells.Item(18, 6).Value = 0.5
$ws.Cells.Item(18, 7).Value = 0.9013820848921223
$ws.Cells.Item(18, 8).Value = 153943008.8498
$ws.Cells.Item(18, 9).Value = -0.3739

# Row 19
$ws.Cells.Item(19, 4).Value = 281.7069
$ws.Cells.Item(19, 5).Value = 622.686
$ws.Cells.Item(19, 6).Value = 133.02
$ws.Cells.Item(19, 7).Value = 0.4524060280783573
$ws.Cells.Item(19, 8).Value = 7786350.2862
$ws.Cells.Item(19, 9).Value = 0.0728

# Row 20
$ws.Cells.Item(20, 3).Value = "Sph"
$ws.Cells.Item(20, 4).Value = 280.2944
$ws.Cells.Item(20, 5).Value = 588.3081
$ws.Cells.Item(20, 6).Value = 79.98
$ws.Cells.Item(20, 7).Value = 0.4764415108342041
$ws.Cells.Item(20, 8).Value = 105176441.6944
$ws.Cells.Item(20, 9).Value = 0.3092

# Row 21
$ws.Cells.Item(21, 3).Value = "Sph"
$ws.Cells.Item(21, 4).Value = 182.5403
$ws.Cells.Item(21, 5).Value = 207.4202
$ws.Cells.Item(21, 6).Value = 9.789999999999999
$ws.Cells.Item(21, 7).Value = 0.880050737584864
$ws.Cells.Item(21, 8).Value = 33183238.3176
$ws.Cells.Item(21, 9).Value = 0.206

# Row 22
$ws.Cells.Item(22, 4).Value = 209.2492
$ws.Cells.Item(22, 5).Value = 3732.2371
$ws.Cells.Item(22, 6).Value = 786.24
$ws.Cells.Item(22, 7).Value = 0.0560653555477491
$ws.Cells.Item(22, 8).Value = 1090703134.3701
$ws.Cells.Item(22, 9).Value = 0.5199

# Row 23
$ws.Cells.Item(23, 4).Value = 257.9469
$ws.Cells.Item(23, 5).Value = 2485.7772
$ws.Cells.Item(23, 6).Value = 855.39
$ws.Cells.Item(23, 7).Value = 0.1037691149472286
$ws.Cells.Item(23, 8).Value = 35928493.2927
$ws.Cells.Item(23, 9).Value = 0.2719

# Row 24
$ws.Cells.Item(24, 3).Value = "Exp"
$ws.Cells.Item(24, 4).Value = 218.4007
$ws.Cells.Item(24, 5).Value = 2771.9139
$ws.Cells.Item(24, 6).Value = 499.76
$ws.Cells.Item(24, 7).Value = 0.07879057859625438
$ws.Cells.Item(24, 8).Value = 10325054493.9236
$ws.Cells.Item(24, 9).Value = -0.3904

# Row 25
$ws.Cells.Item(25, 4).Value = 278.663
$ws.Cells.Item(25, 5).Value = 442138.6957
$ws.Cells.Item(25, 6).Value = 1409.52
$ws.Cells.Item(25, 7).Value = 0.0006302615055187987
$ws.Cells.Item(25, 8).Value = 17415121.0538
$ws.Cells.Item(25, 9).Value = -2.771

# Row 26
$ws.Cells.Item(26, 3).Value = "Sph"
$ws.Cells.Item(26, 4).Value = 193.6977
$ws.Cells.Item(26, 5).Value = 617.5034000000001
$ws.Cells.Item(26, 6).Value = 34.59
$ws.Cells.Item(26, 7).Value = 0.313678758691855
$ws.Cells.Item(26, 8).Value = 192827783781.4646
$ws.Cells.Item(26, 9).Value = -1.9656

# Row 27
$ws.Cells.Item(27, 3).Value = "Gau"
$ws.Cells.Item(27, 4).Value = 144.5362
$ws.Cells.Item(27, 5).Value = 165.0389
$ws.Cells.Item(27, 6).Value = 2.08
$ws.Cells.Item(27, 7).Value = 0.8757705001669304
$ws.Cells.Item(27, 8).Value = 27909113.9696
$ws.Cells.Item(27, 9).Value = -0.08749999999999999

# Row 28
$ws.Cells.Item(28, 3).Value = "Sph"
$ws.Cells.Item(28, 4).Value = 167.4818
$ws.Cells.Item(28, 5).Value = 209.1783
$ws.Cells.Item(28, 6).Value = 0.92
$ws.Cells.Item(28, 7).Value = 0.8006652697722468
$ws.Cells.Item(28, 8).Value = 3661039.193
$ws.Cells.Item(28, 9).Value = 0.4929

# Row 29
$ws.Cells.Item(29, 4).Value = 174.9313
$ws.Cells.Item(29, 5).Value = 184.6152
$ws.Cells.Item(29, 6).Value = 0.18
$ws.Cells.Item(29, 7).Value = 0.947545489212156
$ws.Cells.Item(29, 8).Value = 633979657.6672
$ws.Cells.Item(29, 9).Value = -2.9386

# New rows (30-31) for PA 2019/2020
# Row 30
$ws.Cells.Item(30, 1).Value = "PA"
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = "Exp"
$ws.Cells.Item(30, 4).Value = 159.9252
$ws.Cells.Item(30, 5).Value = 182.2689
$ws.Cells.Item(30, 6).Value = 3.44
$ws.Cells.Item(30, 7).Value = 0.8774135357156376
$ws.Cells.Item(30, 8).Value = 667720.0107
$ws.Cells.Item(30, 9).Value = 0.3915

# Row 31
$ws.Cells.Item(31, 1).Value = "PA"
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = "Gau"
$ws.Cells.Item(31, 4).Value = 183.5904
$ws.Cells.Item(31, 5).Value = 154253.1845
$ws.Cells.Item(31, 6).Value = 765.24
$ws.Cells.Item(31, 7).Value = 0.00119018871859984
$ws.Cells.Item(31, 8).Value = 2565007.9073
$ws.Cells.Item(31, 9).Value = -0.7467
